$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "image_20250807111728_ppp0.jpg"
$ws.Range("I16").Value = "'642,530,686,574"
$ws.Range("I16").Style = "Normal"

$ws.Range("D17").Value = "image_20250807111728_ppp0.jpg"
$ws.Range("I17").Value = "'794,481,830,525"
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = "'0.70"
$ws.Range("J17").Style = "Normal"
